$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row number, new Price (D) text, new Volume(1h) (E) text.
# $null means that column is unchanged for that row (rows 35 untouched entirely;
# rows 36 and 46 keep their original Price and only update Volume(1h)).
$updates = @(
    @{ Row = 2; D = "26.809.08"; E = "  -1.32%  " }
    @{ Row = 3; D = "1.799.73"; E = "  -1.22%  " }
    @{ Row = 4; D = "1.001"; E = "  -0.31%  " }
    @{ Row = 5; D = "309.97"; E = "  -0.68%  " }
    @{ Row = 6; D = "1.001"; E = "  -0.25%  " }
    @{ Row = 7; D = "0.4443"; E = "  +5.14%  " }
    @{ Row = 8; D = "0.3662"; E = "  -0.48%  " }
    @{ Row = 9; D = "0.07314"; E = "  +1.09%  " }
    @{ Row = 10; D = "0.8550"; E = "  +0.07%  " }
    @{ Row = 11; D = "20.59"; E = "  -1.74%  " }
    @{ Row = 12; D = "1.918.51"; E = "  +5.20%  " }
    @{ Row = 13; D = "6.606"; E = "  -1.31%  " }
    @{ Row = 14; D = "92.15"; E = "  +2.85%  " }
    @{ Row = 15; D = "0.07080"; E = "  +0.11%  " }
    @{ Row = 16; D = "5.286"; E = "  -0.10%  " }
    @{ Row = 17; D = "1.002"; E = "  -0.34%  " }
    @{ Row = 18; D = "0.000008701"; E = "  -1.57%  " }
    @{ Row = 19; D = "1.003"; E = "  -0.09%  " }
    @{ Row = 20; D = "14.83"; E = "  -1.16%  " }
    @{ Row = 21; D = "26.849.87"; E = "  -1.47%  " }
    @{ Row = 22; D = "5.143"; E = "  +0.64%  " }
    @{ Row = 23; D = "10.78"; E = "  -0.69%  " }
    @{ Row = 24; D = "1.989"; E = "  +0.46%  " }
    @{ Row = 25; D = "151.89"; E = "  -0.28%  " }
    @{ Row = 26; D = "18.45"; E = "  +0.54%  " }
    @{ Row = 27; D = "2.179"; E = "  -0.87%  " }
    @{ Row = 28; D = "5.188"; E = "  -0.75%  " }
    @{ Row = 29; D = "116.56"; E = "  +0.36%  " }
    @{ Row = 30; D = "0.08795"; E = "  -0.43%  " }
    @{ Row = 31; D = "0.7446"; E = "  -0.32%  " }
    @{ Row = 32; D = "1.165"; E = "  -1.98%  " }
    @{ Row = 33; D = "2.937"; E = "  -0.26%  " }
    @{ Row = 34; D = "4.450"; E = "  +0.28%  " }
    @{ Row = 35; D = $null; E = $null }
    @{ Row = 36; D = $null; E = "  -1.81%  " }
    @{ Row = 37; D = "0.01959"; E = "  -0.22%  " }
    @{ Row = 38; D = "0.05184"; E = "  -0.92%  " }
    @{ Row = 39; D = "0.5306"; E = "  +5.58%  " }
    @{ Row = 40; D = "2.854"; E = "  -0.24%  " }
    @{ Row = 41; D = "7.018"; E = "  -3.41%  " }
    @{ Row = 42; D = "0.1681"; E = "  -0.78%  " }
    @{ Row = 43; D = "0.5123"; E = "  +8.20%  " }
    @{ Row = 44; D = "8.423"; E = "  -2.46%  " }
    @{ Row = 45; D = "10.55"; E = "  -1.03%  " }
    @{ Row = 46; D = $null; E = "  +4.90%  " }
    @{ Row = 47; D = "105.38"; E = "  -0.87%  " }
    @{ Row = 48; D = "0.9995"; E = "  -0.37%  " }
    @{ Row = 49; D = "1.661"; E = "  +0.16%  " }
    @{ Row = 50; D = "0.06312"; E = "  -1.16%  " }
    @{ Row = 51; D = "0.9149"; E = "  +0.43%  " }
)

foreach ($u in $updates) {
    $row = $u.Row
    if ($null -ne $u.D) {
        # Force the Price column to remain plain text (it stores things like
        # "1.001", "0.07314", "26.809.08" which Excel would otherwise
        # auto-convert into numbers / dates on assignment).
        $cell = $ws.Range("D$row")
        $cell.NumberFormat = "@"
        $cell.Value = $u.D
        $cell.Style = "Normal"
    }
    if ($null -ne $u.E) {
        $ws.Range("E$row").Value = $u.E
    }
}
